$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8095196980875699
$ws.Range("C2").Value = 0.2248015017041212
$ws.Range("D2").Value = 0.01752849954137758
$ws.Range("E2").Value = 0.1154778989324541
$ws.Range("F2").Value = 0.4666915567120213
$ws.Range("I2").Value = 0.3178451623716185
$ws.Range("M2").Value = 0.3456511347932931
$ws.Range("O2").Value = 1.490067553058367

$ws.Range("B3").Value = 0.7078771821252303
$ws.Range("C3").Value = 0.1978619620783775
$ws.Range("D3").Value = 0.01565865828340662
$ws.Range("E3").Value = 0.1108315602517393
$ws.Range("F3").Value = 0.4635617682336814
$ws.Range("I3").Value = 0.3245545926818458
$ws.Range("M3").Value = 0.3069021009581476
$ws.Range("O3").Value = 1.493432419653928

$ws.Range("B4").Value = 0.6452677079924172
$ws.Range("C4").Value = 0.1812406018910337
$ws.Range("D4").Value = 0.01450505921918221
$ws.Range("E4").Value = 0.1081095712119904
$ws.Range("F4").Value = 0.4621030246409248
$ws.Range("I4").Value = 0.3290083440213767
$ws.Range("M4").Value = 0.2831425954706717
$ws.Range("O4").Value = 1.497032555072153

$ws.Range("B5").Value = 0.6197049179057501
$ws.Range("C5").Value = 0.1744474129780542
$ws.Range("D5").Value = 0.01403360736041748
$ws.Range("E5").Value = 0.1070329263859264
$ws.Range("F5").Value = 0.4616246226559824
$ws.Range("I5").Value = 0.330907036048556
$ws.Range("M5").Value = 0.273468687830217
$ws.Range("O5").Value = 1.498884071985628

$ws.Range("B6").Value = 0.6154573261825931
$ws.Range("C6").Value = 0.1733182209734139
$ws.Range("D6").Value = 0.01395524256786729
$ws.Range("E6").Value = 0.1068561091016882
$ws.Range("F6").Value = 0.461552182610184
$ws.Range("I6").Value = 0.3312273629607034
$ws.Range("M6").Value = 0.2718628461408628
$ws.Range("O6").Value = 1.499214694783959

$ws.Range("B7").Value = 0.6449231553294794
$ws.Range("C7").Value = 0.1811490664574364
$ws.Range("D7").Value = 0.01449870647849849
$ws.Range("E7").Value = 0.1080949196839569
$ws.Range("F7").Value = 0.462096103359471
$ws.Range("I7").Value = 0.3290336116940438
$ws.Range("M7").Value = 0.2830120960949642
$ws.Range("O7").Value = 1.497055970491886

$ws.Range("B8").Value = 0.7745158556203364
$ws.Range("C8").Value = 0.2155296374945408
$ws.Range("D8").Value = 0.01688494164687881
$ws.Range("E8").Value = 0.1138484804919813
$ws.Range("F8").Value = 0.465516095465766
$ws.Range("I8").Value = 0.3200890637653906
$ws.Range("M8").Value = 0.3322837554056335
$ws.Range("O8").Value = 1.490908555598395

$ws.Range("B9").Value = 1.027001681524382
$ws.Range("C9").Value = 0.2822995273415074
$ws.Range("D9").Value = 0.02151934725368676
$ws.Range("E9").Value = 0.1261849593141307
$ws.Range("F9").Value = 0.4759137672200069
$ws.Range("I9").Value = 0.3052118937926593
$ws.Range("M9").Value = 0.4291645336342214
$ws.Range("O9").Value = 1.491087815790308

$ws.Range("B10").Value = 1.211446582613348
$ws.Range("C10").Value = 0.3309468636901443
$ws.Range("D10").Value = 0.02489542674189238
$ws.Range("E10").Value = 0.1359131060903707
$ws.Range("F10").Value = 0.4858292803131619
$ws.Range("I10").Value = 0.2959211570312767
$ws.Range("M10").Value = 0.5005096146014836
$ws.Range("O10").Value = 1.498766673317988

$ws.Range("B11").Value = 1.295115976528166
$ws.Range("C11").Value = 0.352986894240189
$ws.Range("D11").Value = 0.02642476604653865
$ws.Range("E11").Value = 0.140487820143818
$ws.Range("F11").Value = 0.4908400923093694
$ws.Range("I11").Value = 0.2920539943755429
$ws.Range("M11").Value = 0.5330049504740515
$ws.Range("O11").Value = 1.503917999767623

$ws.Range("B12").Value = 1.326764311613317
$ws.Range("C12").Value = 0.3613196692631959
$ws.Range("D12").Value = 0.02700292933669601
$ws.Range("E12").Value = 0.1422419899096212
$ws.Range("F12").Value = 0.492809907000975
$ws.Range("I12").Value = 0.2906415471166248
$ws.Range("M12").Value = 0.5453158901342903
$ws.Range("O12").Value = 1.506108606981968

$ws.Range("B13").Value = 1.319949877970998
$ws.Range("C13").Value = 0.3595256548809402
$ws.Range("D13").Value = 0.02687845498411434
$ws.Range("E13").Value = 0.141863221795596
$ws.Range("F13").Value = 0.4923824489407025
$ws.Range("I13").Value = 0.2909434269713778
$ws.Range("M13").Value = 0.5426642588557655
$ws.Range("O13").Value = 1.505626127055564

$ws.Range("B14").Value = 1.297720423541421
$ws.Range("C14").Value = 0.3536727054980986
$ws.Range("D14").Value = 0.02647235141959925
$ws.Range("E14").Value = 0.1406316973284447
$ws.Range("F14").Value = 0.4910006982771762
$ws.Range("I14").Value = 0.2919367483220157
$ws.Range("M14").Value = 0.5340176660181726
$ws.Range("O14").Value = 1.504093405761211

$ws.Range("B15").Value = 1.284099577754148
$ws.Range("C15").Value = 0.3500858609158115
$ws.Range("D15").Value = 0.02622347461779384
$ws.Range("E15").Value = 0.1398802057121387
$ws.Range("F15").Value = 0.4901637673437875
$ws.Range("I15").Value = 0.2925519625321265
$ws.Range("M15").Value = 0.5287221124437309
$ws.Range("O15").Value = 1.503185856788889

$ws.Range("B16").Value = 1.20597371616293
$ws.Range("C16").Value = 0.329504654474249
$ws.Range("D16").Value = 0.02479534782791148
$ws.Range("E16").Value = 0.1356171726911981
$ws.Range("F16").Value = 0.485511911534573
$ws.Range("I16").Value = 0.2961811401197032
$ws.Range("M16").Value = 0.4983867658202144
$ws.Range("O16").Value = 1.498463514346895

$ws.Range("B17").Value = 1.157984578316075
$ws.Range("C17").Value = 0.316855456605964
$ws.Range("D17").Value = 0.02391755903910564
$ws.Range("E17").Value = 0.1330404368752838
$ws.Range("F17").Value = 0.4827865469423429
$ws.Range("I17").Value = 0.2984997502115974
$ws.Range("M17").Value = 0.4797872078058418
$ws.Range("O17").Value = 1.495992254748501

$ws.Range("B18").Value = 1.130360420979912
$ws.Range("C18").Value = 0.3095715344622647
$ws.Range("D18").Value = 0.02341207258230327
$ws.Range("E18").Value = 0.1315724012997208
$ws.Range("F18").Value = 0.4812660461845439
$ws.Range("I18").Value = 0.2998671485818107
$ws.Range("M18").Value = 0.4690930144836187
$ws.Range("O18").Value = 1.494726839439039

$ws.Range("B19").Value = 1.121003613107973
$ws.Range("C19").Value = 0.3071038859982025
$ws.Range("D19").Value = 0.02324082064773592
$ws.Range("E19").Value = 0.131077749469803
$ws.Range("F19").Value = 0.4807593005558601
$ws.Range("I19").Value = 0.3003359212280472
$ws.Range("M19").Value = 0.4654727969526391
$ws.Range("O19").Value = 1.494325137825427

$ws.Range("B20").Value = 1.163095400104737
$ws.Range("C20").Value = 0.3182028613953776
$ws.Range("D20").Value = 0.0240110640890947
$ws.Range("E20").Value = 0.133313279722131
$ws.Range("F20").Value = 0.4830717935310389
$ws.Range("I20").Value = 0.2982494303234411
$ws.Range("M20").Value = 0.4817667720185597
$ws.Range("O20").Value = 1.496239169741216

$ws.Range("B21").Value = 1.304250727738463
$ws.Range("C21").Value = 0.3553922227474686
$ws.Range("D21").Value = 0.02659166033579652
$ws.Range("E21").Value = 0.140992830793131
$ws.Range("F21").Value = 0.491404585764343
$ws.Range("I21").Value = 0.2916435730648814
$ws.Range("M21").Value = 0.5365572271058596
$ws.Range("O21").Value = 1.504537080346552

$ws.Range("B22").Value = 1.396296551371961
$ws.Range("C22").Value = 0.3796198878458767
$ws.Range("D22").Value = 0.02827259199222709
$ws.Range("E22").Value = 0.1461392369729495
$ws.Range("F22").Value = 0.4972723205860632
$ws.Range("I22").Value = 0.2876292866372694
$ws.Range("M22").Value = 0.5723989087632191
$ws.Range("O22").Value = 1.511359311726636

$ws.Range("B23").Value = 1.347189452966575
$ws.Range("C23").Value = 0.3666963641755956
$ws.Range("D23").Value = 0.02737597470849806
$ws.Range("E23").Value = 0.1433807319422655
$ws.Range("F23").Value = 0.4941018756241391
$ws.Range("I23").Value = 0.2897439595537641
$ws.Range("M23").Value = 0.5532665473143226
$ws.Range("O23").Value = 1.507589659789772

$ws.Range("B24").Value = 1.160784905117112
$ws.Range("C24").Value = 0.3175937362001378
$ws.Range("D24").Value = 0.02396879305159416
$ws.Range("E24").Value = 0.1331898858715377
$ws.Range("F24").Value = 0.4829426892080164
$ws.Range("I24").Value = 0.2983624928167607
$ws.Range("M24").Value = 0.4808718143371493
$ws.Range("O24").Value = 1.496127055689186

$ws.Range("B25").Value = 0.9588792958921886
$ws.Range("C25").Value = 0.264307386055151
$ws.Range("D25").Value = 0.02027057965521095
$ws.Range("E25").Value = 0.1227325180174432
$ws.Range("F25").Value = 0.4727030551500135
$ws.Range("I25").Value = 0.3089498032192175
$ws.Range("M25").Value = 0.4029270774830707
$ws.Range("O25").Value = 1.489720359830585
